# new restrictions for transferring money
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a numeric-looking value into a cell as plain TEXT (no
# NumberFormat applied to the destination cell, so no new style gets
# introduced -- matches the source file's plain inlineStr cells).
# We stage the text in a scratch cell via TEXT(), copy it, and paste
# only the resulting value into the destination, then wipe the scratch
# cell completely.
function Set-TextValue($cellAddr, $textValue) {
    $scratch = $ws.Range("ZZ500")
    $scratch.Formula = "=TEXT(" + $textValue + ",""0"")"
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.ClearContents()
}

# --- Row 2 updates -----------------------------------------------------
# Event id changed; Stringified Details + Hash regenerated to reflect it.
Set-TextValue "E2" 231220212
$ws.Range("G2").Value = "1/1.0,?:spuckhafte_ferwirklung#7109(231220212)/0000000000000000000000000000000000000000000000000000000000000000/22489572/~~2005cce2777be4eb4c770fb377d4963b8a24b72ae10a43fd97aaeb33f8dfa5cb"
$ws.Range("H2").Value = "2005cce2777be4eb4c770fb377d4963b8a24b72ae10a43fd97aaeb33f8dfa5cb"

# --- Row 3: new ledger entry -------------------------------------------
# A transfer of 1 from spuckhafte_ferwirklung#7109 to Ayano#3463.
$ws.Range("A3").Value = 2
Set-TextValue "B3" 1
$ws.Range("C3").Value = "spuckhafte_ferwirklung#7109"
$ws.Range("D3").Value = "Ayano#3463"
Set-TextValue "E3" 231220212
$ws.Range("F3").Value = "2005cce2777be4eb4c770fb377d4963b8a24b72ae10a43fd97aaeb33f8dfa5cb"
$ws.Range("G3").Value = "2/1,spuckhafte_ferwirklung#7109:Ayano#3463(231220212)/2005cce2777be4eb4c770fb377d4963b8a24b72ae10a43fd97aaeb33f8dfa5cb/15007020/~~2005cc034ce6dd8803c2bf65ad313694059dea96b43a5deda13fdfc8d2d61f6d"
$ws.Range("H3").Value = "2005cc034ce6dd8803c2bf65ad313694059dea96b43a5deda13fdfc8d2d61f6d"

# Column H needs to be widened to fit the new, much longer hash/details text
# (64.3 is the closest COM-settable ColumnWidth to the target 65.21875 char
# width given this engine's pixel-snapping of the stored value)
$ws.Columns.Item(8).ColumnWidth = 64.3

# Reflect the newly populated rows in the active selection
$ws.Range("A2:XFD6").Select()
